# Update "想去人数" (interest count) figures and sold-out status
# for the 展览 (Exhibitions) and 全部类型 (All types) sheets.
# This mirrors a refreshed data pull (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")

# G2: sold-out label changes from "已售罄" (Sold out) to "不可售" (Not for sale)
$ws1.Range("G2").Value = "不可售"

# F column: updated "want to go" counts
$ws1.Range("F3").Value  = 1852   # 南昌·ACG CLUB动漫游戏嘉年华: 1851 -> 1852
$ws1.Range("F6").Value  = 171    # 九江·首届萤火之星国风动漫嘉年华: 170 -> 171
$ws1.Range("F7").Value  = 2507   # 南昌·CM02动漫游戏博览会: 2502 -> 2507
$ws1.Range("F11").Value = 1509   # 南昌·次元之门动漫游戏嘉年华SP：代号序章: 1507 -> 1509
$ws1.Range("F21").Value = 8      # 鹰潭·BM次元盛典运动番only: 7 -> 8
$ws1.Range("F23").Value = 47     # 南昌·漫拥动漫嘉年华Pro-追光启航: 45 -> 47
$ws1.Range("F24").Value = 1589   # 江西·次元星河国风动漫游戏嘉年华: 1581 -> 1589
$ws1.Range("F26").Value = 390    # 南昌·幻梦境国际动漫游戏嘉年华1th: 389 -> 390
$ws1.Range("F28").Value = 202    # 九江·第一届异次元动漫嘉年华: 199 -> 202

# --- Sheet "全部类型" ---
# This sheet contains the same events but with one extra row (row 4) inserted,
# so matching rows from row 4 onward are shifted down by one versus "展览".
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("G2").Value = "不可售"

$ws4.Range("F3").Value  = 1852   # 南昌·ACG CLUB动漫游戏嘉年华: 1851 -> 1852
$ws4.Range("F7").Value  = 171    # 九江·首届萤火之星国风动漫嘉年华: 170 -> 171
$ws4.Range("F8").Value  = 2507   # 南昌·CM02动漫游戏博览会: 2502 -> 2507
$ws4.Range("F12").Value = 1509   # 南昌·次元之门动漫游戏嘉年华SP：代号序章: 1507 -> 1509
$ws4.Range("F22").Value = 8      # 鹰潭·BM次元盛典运动番only: 7 -> 8
$ws4.Range("F24").Value = 47     # 南昌·漫拥动漫嘉年华Pro-追光启航: 45 -> 47
$ws4.Range("F25").Value = 1589   # 江西·次元星河国风动漫游戏嘉年华: 1581 -> 1589
$ws4.Range("F27").Value = 390    # 南昌·幻梦境国际动漫游戏嘉年华1th: 389 -> 390
$ws4.Range("F29").Value = 202    # 九江·第一届异次元动漫嘉年华: 199 -> 202
